$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 146.875
$ws.Range("I2").Value = 127.2
$ws.Range("J2").Value = 179.66667
$ws.Range("K2").Value = 127.2
$ws.Range("L2").Value = 179.66667
$ws.Range("M2").Value = -14.2
$ws.Range("N2").Value = -405.66667
$ws.Range("H12").Value = 233.33333
$ws.Range("I12").Value = 233.33333
$ws.Range("K12").Value = 233.33333
$ws.Range("M12").Value = -63.33332999999999
$ws.Range("H28").Value = 2498.3333
$ws.Range("J28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("N28").ClearContents()
$ws.Range("H33").Value = 241.23077
$ws.Range("I33").Value = 198.72728
$ws.Range("J33").Value = 475
$ws.Range("K33").Value = 198.72728
$ws.Range("L33").Value = 475
$ws.Range("M33").Value = 30.27271999999999
$ws.Range("N33").Value = -933
$ws.Range("H40").Value = 1688.5
$ws.Range("J40").Value = 1688.4286
$ws.Range("L40").Value = 1688.4286
$ws.Range("N40").Value = -2038.4286
$ws.Range("H62").Value = 999
$ws.Range("I62").Value = 999
$ws.Range("K62").Value = 999
$ws.Range("M62").Value = -375
$ws.Range("H65").Value = 999
$ws.Range("I65").Value = 999
$ws.Range("K65").Value = 4995
$ws.Range("M65").Value = -1875
$ws.Range("H107").Value = 415.27777
$ws.Range("I107").Value = 406.46667
$ws.Range("K107").Value = 406.46667
$ws.Range("M107").Value = 1513.53333
$ws.Range("H111").Value = 5507
$ws.Range("I111").Value = 5507
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 16521
$ws.Range("L111").Value = 0
$ws.Range("M111").Value = -13454
$ws.Range("N111").ClearContents()
$ws.Range("H113").Value = 7750.1816
$ws.Range("I113").Value = 6978.143
$ws.Range("J113").Value = 9101.25
$ws.Range("K113").Value = 6978.143
$ws.Range("L113").Value = 9101.25
$ws.Range("M113").Value = -3724.143
$ws.Range("N113").Value = -15609.25
$ws.Range("H115").Value = 991.6
$ws.Range("I115").Value = 991.6
$ws.Range("K115").Value = 2974.8
$ws.Range("M115").Value = -1407.8
$ws.Range("H116").Value = 8666.333000000001
$ws.Range("I116").Value = 9000
$ws.Range("J116").Value = 8499.5
$ws.Range("K116").Value = 9000
$ws.Range("L116").Value = 8499.5
$ws.Range("M116").Value = -5558
$ws.Range("N116").Value = -15383.5
$ws.Range("H132").Value = 918.0833
$ws.Range("I132").Value = 918.0833
$ws.Range("K132").Value = 2754.2499
$ws.Range("M132").Value = -224.2498999999998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 882.3
$ws.Range("I2").Value = 882.3
$ws.Range("K2").Value = 882.3
$ws.Range("M2").Value = -769.3
$ws.Range("H4").Value = 531.25
$ws.Range("I4").Value = 375
$ws.Range("J4").Value = 1000
$ws.Range("K4").Value = 375
$ws.Range("L4").Value = 1000
$ws.Range("M4").Value = -259
$ws.Range("N4").Value = -1232
$ws.Range("H5").Value = 999
$ws.Range("I5").Value = 999
$ws.Range("K5").Value = 999
$ws.Range("M5").Value = -887
$ws.Range("H74").Value = 2515.75
$ws.Range("J74").Value = 750
$ws.Range("L74").Value = 750
$ws.Range("N74").Value = -2498
$ws.Range("H77").Value = 2515.75
$ws.Range("J77").Value = 750
$ws.Range("L77").Value = 750
$ws.Range("N77").Value = -12486
$ws.Range("H116").Value = 882.3
$ws.Range("I116").Value = 882.3
$ws.Range("K116").Value = 882.3
$ws.Range("M116").Value = 1411.7
$ws.Range("H132").Value = 4806
$ws.Range("I132").Value = 4806
$ws.Range("K132").Value = 14418
$ws.Range("M132").Value = -11888

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 882.3
$ws.Range("I3").Value = 882.3
$ws.Range("K3").Value = 882.3
$ws.Range("M3").Value = -768.3
$ws.Range("H4").Value = 999
$ws.Range("I4").Value = 999
$ws.Range("K4").Value = 999
$ws.Range("M4").Value = -884
$ws.Range("H22").Value = 1083.1666
$ws.Range("I22").Value = 1033
$ws.Range("J22").Value = 1133.3334
$ws.Range("K22").Value = 1033
$ws.Range("L22").Value = 1133.3334
$ws.Range("M22").Value = -860
$ws.Range("N22").Value = -1479.3334
$ws.Range("H107").Value = 1755
$ws.Range("I107").Value = 889
$ws.Range("K107").Value = 889
$ws.Range("M107").Value = 1031
$ws.Range("H134").Value = 1935.9333
$ws.Range("I134").Value = 1935.9333
$ws.Range("K134").Value = 5807.7999
$ws.Range("M134").Value = -3272.7999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H33").Value = 49245.332
$ws.Range("J33").Value = 58666.668
$ws.Range("L33").Value = 58666.668
$ws.Range("N33").Value = -59424.668
$ws.Range("H107").Value = 1024.5
$ws.Range("I107").Value = 799.5
$ws.Range("K107").Value = 799.5
$ws.Range("M107").Value = 1120.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 68.8
$ws.Range("I2").Value = 50
$ws.Range("J2").Value = 73.5
$ws.Range("K2").Value = 300
$ws.Range("L2").Value = 441
$ws.Range("M2").Value = -187
$ws.Range("N2").Value = -667
$ws.Range("H12").Value = 167.33333
$ws.Range("I12").Value = 49
$ws.Range("J12").Value = 178.09091
$ws.Range("K12").Value = 147
$ws.Range("L12").Value = 534.27273
$ws.Range("M12").Value = 26
$ws.Range("N12").Value = -880.27273
$ws.Range("H14").Value = 2100.1538
$ws.Range("I14").Value = 2100.1538
$ws.Range("K14").Value = 6300.4614
$ws.Range("M14").Value = -6127.4614

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 193.27272
$ws.Range("I2").Value = 231.22223
$ws.Range("K2").Value = 231.22223
$ws.Range("M2").Value = -118.22223
$ws.Range("H70").Value = 3000
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 3000
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 3000
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -3540
$ws.Range("H73").Value = 3000
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 3000
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 3000
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = -4872
$ws.Range("H113").Value = 4999.25
$ws.Range("I113").Value = 4999
$ws.Range("J113").Value = 4999.5
$ws.Range("K113").Value = 4999
$ws.Range("L113").Value = 4999.5
$ws.Range("M113").Value = -2829
$ws.Range("N113").Value = -9339.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 999.4
$ws.Range("J22").Value = 1000
$ws.Range("L22").Value = 1000
$ws.Range("N22").Value = -1590
$ws.Range("H27").Value = 999.4
$ws.Range("J27").Value = 1000
$ws.Range("L27").Value = 1000
$ws.Range("N27").Value = -1214
$ws.Range("H61").Value = 1223.25
$ws.Range("I61").Value = 1004
$ws.Range("J61").Value = 1442.5
$ws.Range("K61").Value = 1004
$ws.Range("L61").Value = 1442.5
$ws.Range("M61").Value = -802
$ws.Range("N61").Value = -1846.5
$ws.Range("H113").Value = 1223.25
$ws.Range("I113").Value = 1004
$ws.Range("J113").Value = 1442.5
$ws.Range("K113").Value = 1004
$ws.Range("L113").Value = 1442.5
$ws.Range("M113").Value = 1166
$ws.Range("N113").Value = -5782.5
$ws.Range("H132").Value = 2071.4285
$ws.Range("I132").Value = 1916.6666
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 5749.9998
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -3219.9998
$ws.Range("N132").Value = -14060

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1439.5714
$ws.Range("I107").Value = 1362.8334
$ws.Range("J107").Value = 1900
$ws.Range("K107").Value = 4088.5002
$ws.Range("L107").Value = 5700
$ws.Range("M107").Value = -2168.5002
$ws.Range("N107").Value = -9540
$ws.Range("H113").Value = 273.625
$ws.Range("I113").Value = 266.5
$ws.Range("J113").Value = 295
$ws.Range("K113").Value = 799.5
$ws.Range("L113").Value = 885
$ws.Range("M113").Value = 1370.5
$ws.Range("N113").Value = -5225
$ws.Range("H124").Value = 60000
$ws.Range("J124").Value = 50000
$ws.Range("L124").Value = 50000
$ws.Range("N124").Value = -59820
$ws.Range("H136").Value = 1856.25
$ws.Range("I136").Value = 1752.2727
$ws.Range("K136").Value = 5256.8181
$ws.Range("M136").Value = -2706.8181
